$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 10.86
$ws.Range("D3").Value = 9.039999999999999
$ws.Range("D4").Value = 12.41
$ws.Range("D5").Value = 15.42
$ws.Range("D6").Value = 3.06
$ws.Range("D7").Value = 8.539999999999999
$ws.Range("D8").Value = 7.66
$ws.Range("D9").Value = 5.64
$ws.Range("D10").Value = 25.09
$ws.Range("D11").Value = 1.23
$ws.Range("D12").Value = 3.14
$ws.Range("D13").Value = 1.25
$ws.Range("D14").Value = 6.46
$ws.Range("D15").Value = 3.8
$ws.Range("D16").Value = 4.68
$ws.Range("D17").Value = 16.1
$ws.Range("D18").Value = 3.38
$ws.Range("D19").Value = 13.12
$ws.Range("D20").Value = 4.85
$ws.Range("D21").Value = 1.35
$ws.Range("D22").Value = 7.26
$ws.Range("D23").Value = 3.02
$ws.Range("D24").Value = 3.4
$ws.Range("D26").Value = 5.19
$ws.Range("D28").Value = 4.57
$ws.Range("D29").Value = 0.73
$ws.Range("D30").Value = 1.2
$ws.Range("D31").Value = 2.13
$ws.Range("D32").Value = 1.44
$ws.Range("D33").Value = 3.07
$ws.Range("D34").Value = 11.02
$ws.Range("D36").Value = 1.89
$ws.Range("D37").Value = 5.38
$ws.Range("D38").Value = 26.58
$ws.Range("D39").Value = 1.8
$ws.Range("D40").Value = 1.01
$ws.Range("D43").Value = 3.31
$ws.Range("D44").Value = 0.87
$ws.Range("D45").Value = 1.12
$ws.Range("D46").Value = 3.41
$ws.Range("D47").Value = 7.9
$ws.Range("D48").Value = 4.22
$ws.Range("D49").Value = 1.3
$ws.Range("D50").Value = 1.9
$ws.Range("D51").Value = 6.09
$ws.Range("D52").Value = 1.96
$ws.Range("D53").Value = 3.62
$ws.Range("D54").Value = 4.04
$ws.Range("D55").Value = 1.38
$ws.Range("D56").Value = 10.86
$ws.Range("D57").Value = 8.07
$ws.Range("D58").Value = 16.79
$ws.Range("D59").Value = 1.66
$ws.Range("D60").Value = 9.16
$ws.Range("D61").Value = 8.94
$ws.Range("D62").Value = 3.71
$ws.Range("D63").Value = 1.9
$ws.Range("D64").Value = 6.06
$ws.Range("D65").Value = 3.53
$ws.Range("D66").Value = 1.02
$ws.Range("D67").Value = 0.02
